$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 4925.75
$ws.Range("I20").Value = 772.2857
$ws.Range("K20").Value = 772.2857
$ws.Range("M20").Value = -542.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 4925.75
$ws.Range("I35").Value = 772.2857
$ws.Range("K35").Value = 772.2857
$ws.Range("M35").Value = -393.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 32400.75
$ws.Range("J108").Value = 32400.75
$ws.Range("L108").Value = 32400.75
$ws.Range("N108").Value = -40080.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H120").Value = 49365
$ws.Range("J120").Value = 49365
$ws.Range("L120").Value = 49365
$ws.Range("N120").Value = -59041

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 13763.013
$ws.Range("I132").Value = 2377.1045
$ws.Range("J132").Value = 98524.78
$ws.Range("K132").Value = 7131.3135
$ws.Range("L132").Value = 295574.34
$ws.Range("M132").Value = -4601.3135
$ws.Range("N132").Value = -300634.34

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22126.492
$ws.Range("I32").Value = 20155.188
$ws.Range("K32").Value = 20155.188
$ws.Range("M32").Value = -19868.188

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1897.234
$ws.Range("I74").Value = 1445.5278
$ws.Range("J74").Value = 3375.5454
$ws.Range("K74").Value = 1445.5278
$ws.Range("L74").Value = 3375.5454
$ws.Range("M74").Value = -571.5278000000001
$ws.Range("N74").Value = -5123.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1897.234
$ws.Range("I77").Value = 1445.5278
$ws.Range("J77").Value = 3375.5454
$ws.Range("K77").Value = 7227.639
$ws.Range("L77").Value = 16877.727
$ws.Range("M77").Value = -2859.639
$ws.Range("N77").Value = -25613.727

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9092620
$ws.Range("I132").Value = 13514603
$ws.Range("J132").Value = 2986.5557
$ws.Range("K132").Value = 40543809
$ws.Range("L132").Value = 8959.667099999999
$ws.Range("M132").Value = -40541279
$ws.Range("N132").Value = -14019.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 46711
$ws.Range("J119").Value = 46711
$ws.Range("L119").Value = 46711
$ws.Range("N119").Value = -56387

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H120").Value = 45761
$ws.Range("J120").Value = 45761
$ws.Range("L120").Value = 45761
$ws.Range("N120").Value = -55437

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5199.359
$ws.Range("I31").Value = 2558.44
$ws.Range("J31").Value = 9915.286
$ws.Range("K31").Value = 2558.44
$ws.Range("L31").Value = 9915.286
$ws.Range("M31").Value = -2263.44
$ws.Range("N31").Value = -10505.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5199.359
$ws.Range("I34").Value = 2558.44
$ws.Range("J34").Value = 9915.286
$ws.Range("K34").Value = 2558.44
$ws.Range("L34").Value = 9915.286
$ws.Range("M34").Value = -2356.44
$ws.Range("N34").Value = -10319.286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 49724.35
$ws.Range("I132").Value = 12505.833
$ws.Range("J132").Value = 91595.19
$ws.Range("K132").Value = 37517.499
$ws.Range("L132").Value = 274785.57
$ws.Range("M132").Value = -34987.499
$ws.Range("N132").Value = -279845.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 625.375
$ws.Range("I47").Value = 533.8333
$ws.Range("K47").Value = 1601.4999
$ws.Range("M47").Value = -1170.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 2500
$ws.Range("J93").Value = 2500
$ws.Range("L93").Value = 7500
$ws.Range("N93").Value = -11244

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7113.7095
$ws.Range("I107").Value = 13047.875
$ws.Range("J107").Value = 5049.6523
$ws.Range("K107").Value = 39143.625
$ws.Range("L107").Value = 15148.9569
$ws.Range("M107").Value = -37223.625
$ws.Range("N107").Value = -18988.9569

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 17666.666
$ws.Range("I113").Value = 26050
$ws.Range("J113").Value = 900
$ws.Range("K113").Value = 78150
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -75980
$ws.Range("N113").Value = -7040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1174.1428
$ws.Range("J117").Value = 797.25
$ws.Range("L117").Value = 2391.75
$ws.Range("N117").Value = -9275.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 103923.266
$ws.Range("I121").Value = 303.75
$ws.Range("J121").Value = 141603.1
$ws.Range("K121").Value = 911.25
$ws.Range("L121").Value = 424809.3
$ws.Range("M121").Value = 398.75
$ws.Range("N121").Value = -427429.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 905.1
$ws.Range("J131").Value = 935.1183
$ws.Range("L131").Value = 2805.3549
$ws.Range("N131").Value = -12885.3549

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1453.6842
$ws.Range("I132").Value = 829.4
$ws.Range("J132").Value = 1860.826
$ws.Range("K132").Value = 7464.599999999999
$ws.Range("L132").Value = 16747.434
$ws.Range("M132").Value = -4934.599999999999
$ws.Range("N132").Value = -21807.434

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4030.9736
$ws.Range("I132").Value = 4492.476
$ws.Range("J132").Value = 3460.8823
$ws.Range("K132").Value = 13477.428
$ws.Range("L132").Value = 10382.6469
$ws.Range("M132").Value = -10947.428
$ws.Range("N132").Value = -15442.6469

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1546.9524
$ws.Range("I100").Value = 1339.0667
$ws.Range("K100").Value = 1339.0667
$ws.Range("M100").Value = -798.0667000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 42417
$ws.Range("J123").Value = 42417
$ws.Range("L123").Value = 42417
$ws.Range("N123").Value = -52217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4831.6333
$ws.Range("I132").Value = 4778
$ws.Range("J132").Value = 4938.9
$ws.Range("K132").Value = 14334
$ws.Range("L132").Value = 14816.7
$ws.Range("M132").Value = -11804
$ws.Range("N132").Value = -19876.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2748.9722
$ws.Range("I136").Value = 2266.7083
$ws.Range("J136").Value = 3713.5
$ws.Range("K136").Value = 6800.124899999999
$ws.Range("L136").Value = 11140.5
$ws.Range("M136").Value = -4250.124899999999
$ws.Range("N136").Value = -16240.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 26090.334
$ws.Range("J69").Value = 26090.334
$ws.Range("L69").Value = 26090.334
$ws.Range("N69").Value = -27588.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H72").Value = 26090.334
$ws.Range("J72").Value = 26090.334
$ws.Range("L72").Value = 78271.00199999999
$ws.Range("N72").Value = -85759.00199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1578.85
$ws.Range("I132").Value = 725.53845
$ws.Range("J132").Value = 3163.5715
$ws.Range("K132").Value = 2176.61535
$ws.Range("L132").Value = 9490.7145
$ws.Range("M132").Value = 353.38465
$ws.Range("N132").Value = -14550.7145
